# Generate Report for Archive
# - Flip the "Ready for handoff" status to "In Translation" everywhere it
#   appears (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 all share the same
#   status string), then shrink the now-narrower status columns to match.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status ---
$wsOverview = $wb.Worksheets.Item("Overview")
$statusRangeOverview = $wsOverview.Range("E2:F4")
$statusRangeOverview.Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: column C holds the status ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: column C holds the status ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
